# Applies targeted cell-value updates to match the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(98, 8).Value = 511.4
$ws.Cells.Item(98, 9).Value = 476.5
$ws.Cells.Item(98, 11).Value = 476.5
$ws.Cells.Item(98, 13).Value = 1021.5
$ws.Cells.Item(106, 8).Value = 3216.1667
$ws.Cells.Item(106, 10).Value = 0
$ws.Cells.Item(106, 12).Value = 0
$ws.Cells.Item(106, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 3239.6
$ws.Cells.Item(113, 9).Value = 3066.3333
$ws.Cells.Item(113, 10).Value = 3499.5
$ws.Cells.Item(113, 11).Value = 3066.3333
$ws.Cells.Item(113, 12).Value = 3499.5
$ws.Cells.Item(113, 13).Value = 187.6667000000002
$ws.Cells.Item(113, 14).Value = -10007.5
$ws.Cells.Item(122, 8).Value = 511.4
$ws.Cells.Item(122, 9).Value = 476.5
$ws.Cells.Item(122, 11).Value = 1429.5
$ws.Cells.Item(122, 13).Value = 1020.5
$ws.Cells.Item(137, 8).Value = 2678.7307
$ws.Cells.Item(137, 9).Value = 1705.6923
$ws.Cells.Item(137, 11).Value = 5117.0769
$ws.Cells.Item(137, 13).Value = -2567.0769

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(61, 8).Value = 3735
$ws.Cells.Item(61, 9).Value = 3735
$ws.Cells.Item(61, 11).Value = 3735
$ws.Cells.Item(61, 13).Value = -3523
$ws.Cells.Item(74, 8).Value = 4452.4
$ws.Cells.Item(74, 9).Value = 4452.4
$ws.Cells.Item(74, 11).Value = 4452.4
$ws.Cells.Item(74, 13).Value = -3578.4
$ws.Cells.Item(77, 8).Value = 4452.4
$ws.Cells.Item(77, 9).Value = 4452.4
$ws.Cells.Item(77, 11).Value = 22262
$ws.Cells.Item(77, 13).Value = -17894
$ws.Cells.Item(97, 8).Value = 1477.3846
$ws.Cells.Item(97, 9).Value = 999.5454999999999
$ws.Cells.Item(97, 10).Value = 4105.5
$ws.Cells.Item(97, 11).Value = 999.5454999999999
$ws.Cells.Item(97, 12).Value = 4105.5
$ws.Cells.Item(97, 13).Value = -503.5454999999999
$ws.Cells.Item(97, 14).Value = -5097.5
$ws.Cells.Item(131, 8).Value = 69993.336
$ws.Cells.Item(131, 10).Value = 69993.336
$ws.Cells.Item(131, 12).Value = 69993.336
$ws.Cells.Item(131, 14).Value = -80073.336
$ws.Cells.Item(136, 8).Value = 3735
$ws.Cells.Item(136, 9).Value = 3735
$ws.Cells.Item(136, 11).Value = 11205
$ws.Cells.Item(136, 13).Value = -8655

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 3008
$ws.Cells.Item(20, 9).Value = 3008
$ws.Cells.Item(20, 11).Value = 3008
$ws.Cells.Item(20, 13).Value = -2761
$ws.Cells.Item(94, 8).Value = 2102.6667
$ws.Cells.Item(94, 9).Value = 2102.6667
$ws.Cells.Item(94, 11).Value = 2102.6667
$ws.Cells.Item(94, 13).Value = -1651.6667

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 6007
$ws.Cells.Item(31, 9).Value = 3132
$ws.Cells.Item(31, 11).Value = 3132
$ws.Cells.Item(31, 13).Value = -2837
$ws.Cells.Item(34, 8).Value = 6007
$ws.Cells.Item(34, 9).Value = 3132
$ws.Cells.Item(34, 11).Value = 3132
$ws.Cells.Item(34, 13).Value = -2930
$ws.Cells.Item(70, 8).Value = 0
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(80, 8).Value = 20116
$ws.Cells.Item(80, 9).Value = 20116
$ws.Cells.Item(80, 11).Value = 20116
$ws.Cells.Item(80, 13).Value = -18993
$ws.Cells.Item(83, 8).Value = 20116
$ws.Cells.Item(83, 9).Value = 20116
$ws.Cells.Item(83, 11).Value = 60348
$ws.Cells.Item(83, 13).Value = -54732
$ws.Cells.Item(99, 8).Value = 2187.0715
$ws.Cells.Item(99, 9).Value = 1406.3334
$ws.Cells.Item(99, 11).Value = 1406.3334
$ws.Cells.Item(99, 13).Value = 91.66660000000002
$ws.Cells.Item(126, 8).Value = 2187.0715
$ws.Cells.Item(126, 9).Value = 1406.3334
$ws.Cells.Item(126, 11).Value = 4219.0002
$ws.Cells.Item(126, 13).Value = -1749.0002

$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(7, 8).Value = 555
$ws.Cells.Item(7, 9).Value = 415
$ws.Cells.Item(7, 10).Value = 695
$ws.Cells.Item(7, 11).Value = 1245
$ws.Cells.Item(7, 12).Value = 2085
$ws.Cells.Item(7, 13).Value = -1133
$ws.Cells.Item(7, 14).Value = -2309
$ws.Cells.Item(16, 8).Value = 10000
$ws.Cells.Item(16, 10).Value = 10000
$ws.Cells.Item(16, 12).Value = 30000
$ws.Cells.Item(16, 14).Value = -30346
$ws.Cells.Item(121, 8).Value = 199
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 11500
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1087.7778
$ws.Cells.Item(107, 9).Value = 1022.5833
$ws.Cells.Item(107, 10).Value = 1218.1666
$ws.Cells.Item(107, 11).Value = 1022.5833
$ws.Cells.Item(107, 12).Value = 1218.1666
$ws.Cells.Item(107, 13).Value = 897.4167
$ws.Cells.Item(107, 14).Value = -5058.1666

$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 2101.75
$ws.Cells.Item(22, 9).Value = 1679.8
$ws.Cells.Item(22, 10).Value = 2805
$ws.Cells.Item(22, 11).Value = 1679.8
$ws.Cells.Item(22, 12).Value = 2805
$ws.Cells.Item(22, 13).Value = -1384.8
$ws.Cells.Item(22, 14).Value = -3395
$ws.Cells.Item(27, 8).Value = 2101.75
$ws.Cells.Item(27, 9).Value = 1679.8
$ws.Cells.Item(27, 10).Value = 2805
$ws.Cells.Item(27, 11).Value = 1679.8
$ws.Cells.Item(27, 12).Value = 2805
$ws.Cells.Item(27, 13).Value = -1572.8
$ws.Cells.Item(27, 14).Value = -3019
$ws.Cells.Item(40, 8).Value = 2345
$ws.Cells.Item(40, 9).Value = 2345
$ws.Cells.Item(40, 11).Value = 2345
$ws.Cells.Item(40, 13).Value = -2209
$ws.Cells.Item(74, 8).Value = 87598.5
$ws.Cells.Item(74, 9).Value = 87598.5
$ws.Cells.Item(74, 11).Value = 87598.5
$ws.Cells.Item(74, 13).Value = -86600.5
$ws.Cells.Item(77, 8).Value = 87598.5
$ws.Cells.Item(77, 9).Value = 87598.5
$ws.Cells.Item(77, 11).Value = 262795.5
$ws.Cells.Item(77, 13).Value = -257803.5

$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(54, 8).Value = 55249.5
$ws.Cells.Item(54, 9).Value = 1000
$ws.Cells.Item(54, 10).Value = 73332.664
$ws.Cells.Item(54, 11).Value = 1000
$ws.Cells.Item(54, 12).Value = 73332.664
$ws.Cells.Item(54, 13).Value = -480
$ws.Cells.Item(54, 14).Value = -74372.664
$ws.Cells.Item(62, 8).Value = 11687.5
$ws.Cells.Item(62, 10).Value = 11714.286
$ws.Cells.Item(62, 12).Value = 11714.286
$ws.Cells.Item(62, 14).Value = -12962.286
$ws.Cells.Item(65, 8).Value = 11687.5
$ws.Cells.Item(65, 10).Value = 11714.286
$ws.Cells.Item(65, 12).Value = 58571.43
$ws.Cells.Item(65, 14).Value = -64811.43
$ws.Cells.Item(75, 8).Value = 60000
$ws.Cells.Item(75, 9).Value = 60000
$ws.Cells.Item(75, 11).Value = 60000
$ws.Cells.Item(75, 13).Value = -59064
$ws.Cells.Item(78, 8).Value = 60000
$ws.Cells.Item(78, 9).Value = 60000
$ws.Cells.Item(78, 11).Value = 180000
$ws.Cells.Item(78, 13).Value = -175320
$ws.Cells.Item(105, 8).Value = 27500
$ws.Cells.Item(105, 10).Value = 27500
$ws.Cells.Item(105, 12).Value = 27500
$ws.Cells.Item(105, 14).Value = -34488
$ws.Cells.Item(125, 8).Value = 31034.25
$ws.Cells.Item(125, 10).Value = 31034.25
$ws.Cells.Item(125, 12).Value = 31034.25
$ws.Cells.Item(125, 14).Value = -40874.25
$ws.Cells.Item(132, 8).Value = 2583.2222
$ws.Cells.Item(132, 9).Value = 2713.3076
$ws.Cells.Item(132, 10).Value = 2245
$ws.Cells.Item(132, 11).Value = 8139.9228
$ws.Cells.Item(132, 12).Value = 6735
$ws.Cells.Item(132, 13).Value = -5609.9228
$ws.Cells.Item(132, 14).Value = -11795
